$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds text values (e.g. "0") rather than numbers, so force text
# formatting before assigning so Excel stores them as strings, not numbers.
$ws.Range("B2:B7").NumberFormat = "@"

$ws.Range("B2").Value = "0"
$ws.Range("B3").Value = "0"
$ws.Range("B4").Value = "0"
$ws.Range("B5").Value = "0"
$ws.Range("B6").Value = "0"
$ws.Range("B7").Value = "0"

$ws.Range("D3").Value = 47
